# BIS-769: Fixed xls test files
# Add two new header columns ("Pattern", "Pattern Type") to the right of the
# existing "Unique" header column (K=Multivalued, L=Unique) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell M4 ("Pattern") - copy L4's format (bold header style) then set its value.
$ws.Cells.Item(4, 12).Copy()
$ws.Cells.Item(4, 13).PasteSpecial(-4122)
$ws.Cells.Item(4, 13).Value = "Pattern"

# New header cell N4 ("Pattern Type") - copy L4's format (bold header style) then set its value.
$ws.Cells.Item(4, 12).Copy()
$ws.Cells.Item(4, 14).PasteSpecial(-4122)
$ws.Cells.Item(4, 14).Value = "Pattern Type"

# Match the selection left behind by the edit (new columns M:N selected).
[void]$ws.Range("M4:N4").Select()
